$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Echoes of Forgotten Verses: Inviting the Muse" "Unveiling the Invisible: The Enigmatic Realm of Physics"
Replace-Text "Allison Monroe" "Henry Morgan"
Replace-Text "allisonm7k@gmail" "henrymorgan8472@gmail"

Replace-Text "The halls of literature whisper with countless tales and songs unheard, remnants of ancient bards and poets echoing amidst the sands of time" "Physics, an enigmatic subject yet a vital force in our lives, orchestrates the dance of atoms, forms the very fabric of reality, and governs the cosmic ballet"

Replace-Text " Poetry, an art form transcending cultures and eras, conjures ethereal worlds and paints emotions with words" " It is often dubbed the language of the universe, for within its equations lies a deeper understanding of how our world operates"

Replace-Text " Its origins trace back to the dawn of human history, weaving narratives that touched hearts and shaped destinies. Like a mystic thread, poetry binds us to our ancestors, weaving a rich tapestry of emotions and shared experiences. Its power lies in its ability to unlock the gateways of imagination and offer solace in the face of life's challenges" " We may perceive the world through our senses, but physics delves into the intricate mechanisms that underpin the phenomena we behold"

Replace-Text "The written word, a timeless conduit for human thought and sentiment, captures the essence of emotions in a manner that transcends spoken language" "In the realm of physics, the intricate symphony of forces, particles, and energy weaves a tapestry of complexity"

Replace-Text " Verse possesses an enrapturing quality, distilling profound thoughts and experiences into succinct, impactful phrases" " From the subatomic realm of quarks and leptons to the vast expanse of galaxies and black holes, physics offers a glimpse into the enigmas that have captivated mankind for ages"

Replace-Text " It immortalizes moments, etching them onto the pages of history and the chambers of memory. From the love sonnets of Petrarch to the poignant reflections of Emily Dickinson, poetry has served as a vessel for the soul, embracing the full spectrum of human emotions, from joy to sorrow, love to loss" " It is a journey to unravel the secrets of the universe, to understand the fundamental laws that shape our existence, and to decipher the cosmic code that governs the symphony of life"

Replace-Text "Throughout history, prominent poets have emerged as beacons of inspiration, their verses resonating with generations" "The study of physics is not merely about equations and formulas; it is an odyssey to question the fabric of reality, to ponder the nature of time, space, and matter"

Replace-Text " From Homer, whose epic sagas brought to life gods and heroes, to Rumi, whose words transcended cultural and religious divides, poetry continues to echo in the hearts of millions. Their artistry invites us to ponder upon the mysteries of existence, challenges our perceptions, and evokes a profound connection with the human condition. In a world rife with distractions, poetry serves as a sanctuary for reflection, a moment of pause amid the incessant noise" " It challenges us to push the boundaries of human knowledge, to explore the mysteries that lie beyond the tangible and visible, and to unveil the hidden patterns that connect the disparate elements of our cosmos"

Replace-Text "Poetry, the art of verse, has spanned time and culture, echoing the human experience in profound and evocative ways" "This essay explored the enigmatic realm of physics, delving into its fundamental concepts and highlighting its significance in comprehending the world around us"

Replace-Text " Rooted in ancient traditions, it transcends spoken language with its ability to capture emotions and thoughts in a concise and impactful manner" " It emphasizes the intricate nature of physics and its unique ability to unravel the secrets of the universe, from the subatomic realm to the vast expanse of galaxies"

Replace-Text " Renowned poets illuminated our understanding of human nature through verse, offering solace, exploration of the human condition, and a chance to pause amidst the world's chaos. Poetry invites us to delve into our imaginations, urging us to reflect on the complexities of life and craft a world that is both introspective and transformative" " Physics poses questions that challenge human knowledge and invites us on a journey to understand the laws that govern existence, pushing the boundaries of human comprehension and exploration"

# Add a new empty paragraph at the end of the document body (before sectPr)
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Host "Edit complete"
